# ajuste de endpoints para que funcionen en el frontend
# Adds three new course/schedule blocks (rows 12-17) to the schedule table,
# following the same 2-row-per-course pattern (merged A/B/C columns, two
# stacked schedule lines in column D) used by the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New course blocks: (topRow, bottomRow, Materia, Nombre, CodigoDocente, Horario1, Horario2)
$courses = @(
    @(12, 13, "1155105A", "INTRODUCCION ING SISTEMAS", 1827, "MARTES 08:00-09:00 SF404", "JUEVES 08:00-10:00 SA403"),
    @(14, 15, "1155201A", "CALCULO INTEGRAL",          7491, "MARTES 06:00-08:00 SA203", "JUEVES 06:00-08:00 SA202"),
    @(16, 17, "1155102A", "MATEMATICAS DISCRETAS",     4412, "MARTES 10:00-12:00 SA402", "MIERCOLES 09:00-10:00 SA414")
)

foreach ($course in $courses) {
    $top    = $course[0]
    $bottom = $course[1]
    $materia = $course[2]
    $nombre  = $course[3]
    $codigo  = $course[4]
    $horario1 = $course[5]
    $horario2 = $course[6]

    # Copy formatting (borders + wrap text) from the last existing course
    # block (rows 10:11) so the new rows look identical to the rest of the
    # table.
    $srcRange = $ws.Range("A10:D11")
    $srcRange.Copy() | Out-Null
    $dstRange = $ws.Range("A" + $top + ":D" + $bottom)
    $dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # Fill in the values.
    $ws.Range("A" + $top).Value = $materia
    $ws.Range("B" + $top).Value = $nombre
    $ws.Range("C" + $top).Value = $codigo
    $ws.Range("D" + $top).Value = $horario1
    $ws.Range("D" + $bottom).Value = $horario2

    # Merge the Materia/Nombre/Codigo columns across the two rows.
    $ws.Range("A" + $top + ":A" + $bottom).Merge() | Out-Null
    $ws.Range("B" + $top + ":B" + $bottom).Merge() | Out-Null
    $ws.Range("C" + $top + ":C" + $bottom).Merge() | Out-Null
}

# Set the page to portrait orientation.
$ws.PageSetup.Orientation = 1

# Restore the selection to where the user last left it.
$ws.Range("D26").Select() | Out-Null

Write-Output "done"
